$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date for 017d85ba... row (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 16:49:36"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 017d85ba... row (row 3)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-23 16:49:31"
$wsZhCn.Range("K3").Value = "2016-08-23 16:49:49"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime for 017d85ba... row (row 3)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-23 16:49:36"
$wsDeDe.Range("K3").Value = "2016-08-23 16:49:57"
